$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A8").Value = 46070
$ws.Range("D8").Value = 158.33000000000001
$ws.Range("E8").Value = 149.19
$ws.Range("F8").Value = 159.19
$ws.Range("G8").Value = 149.08000000000001

$ws.Range("A9").Value = 46070
$ws.Range("D9").Value = 158.33000000000001
$ws.Range("E9").Value = 149.19
$ws.Range("F9").Value = 159.19
$ws.Range("G9").Value = 149.08000000000001

$ws.Range("A10").Value = 46070
$ws.Range("D10").Value = 159.71
$ws.Range("E10").Value = 151.87
$ws.Range("F10").Value = 161.87
$ws.Range("G10").Value = 152.11000000000001

$ws.Range("A11").Value = 46067
$ws.Range("D11").Value = 159
$ws.Range("E11").Value = 149.91
$ws.Range("F11").Value = 159.91
$ws.Range("G11").Value = 149.80000000000001

$ws.Range("A12").Value = 46067
$ws.Range("D12").Value = 159
$ws.Range("E12").Value = 149.91
$ws.Range("F12").Value = 159.91
$ws.Range("G12").Value = 149.80000000000001

$ws.Range("A13").Value = 46067
$ws.Range("D13").Value = 160.56
$ws.Range("E13").Value = 152.55000000000001
$ws.Range("F13").Value = 162.55000000000001
$ws.Range("G13").Value = 152.79

$ws.Range("A17").Value = 46070
$ws.Range("D17").Value = 164.07
$ws.Range("E17").Value = 155.41
$ws.Range("F17").Value = 165.41

$ws.Range("A18").Value = 46067
$ws.Range("D18").Value = 164.89
$ws.Range("E18").Value = 156.06
$ws.Range("F18").Value = 166.06

$ws.Range("A22").Value = 46070
$ws.Range("D22").Value = 159.51
$ws.Range("E22").Value = 151.53
$ws.Range("F22").Value = 161.13
$ws.Range("G22").Value = 153.29

$ws.Range("A23").Value = 46070
$ws.Range("D23").Value = 164.68
$ws.Range("E23").Value = 157.57
$ws.Range("F23").Value = 167.57

$ws.Range("A24").Value = 46070
$ws.Range("D24").Value = 164.87
$ws.Range("E24").Value = 158.11000000000001
$ws.Range("F24").Value = 168.11

$ws.Range("A25").Value = 46070
$ws.Range("D25").Value = 164.88
$ws.Range("E25").Value = 157.63
$ws.Range("F25").Value = 167.63
$ws.Range("G25").Value = 158.47999999999999

$ws.Range("A26").Value = 46070
$ws.Range("D26").Value = 164.5
$ws.Range("E26").Value = 159.22
$ws.Range("F26").Value = 169.22

$ws.Range("A27").Value = 46067
$ws.Range("D27").Value = 160.18
$ws.Range("E27").Value = 152.03
$ws.Range("F27").Value = 161.63
$ws.Range("G27").Value = 153.78

$ws.Range("A28").Value = 46067
$ws.Range("D28").Value = 165.54
$ws.Range("E28").Value = 158.25
$ws.Range("F28").Value = 168.25

$ws.Range("A29").Value = 46067
$ws.Range("D29").Value = 165.72

$ws.Range("A30").Value = 46067
$ws.Range("D30").Value = 165.72

$ws.Range("A31").Value = 46067
$ws.Range("D31").Value = 165.35

$ws.Range("A35").Value = 46070
$ws.Range("D35").Value = 158.22
$ws.Range("E35").Value = 149.56
$ws.Range("F35").Value = 158.56

$ws.Range("A36").Value = 46067
$ws.Range("E36").Value = 150.24
$ws.Range("F36").Value = 159.24

$ws.Range("A40").Value = 46070
$ws.Range("D40").Value = 164.42
$ws.Range("E40").Value = 156.85
$ws.Range("F40").Value = 166.85

$ws.Range("A41").Value = 46070
$ws.Range("D41").Value = 164.14
$ws.Range("E41").Value = 157.27000000000001
$ws.Range("F41").Value = 167.27

$ws.Range("A42").Value = 46067
$ws.Range("D42").Value = 165.3
$ws.Range("E42").Value = 157.66
$ws.Range("F42").Value = 167.66

$ws.Range("A43").Value = 46067
$ws.Range("D43").Value = 165.02
$ws.Range("E43").Value = 158.08000000000001
$ws.Range("F43").Value = 168.08

$ws.Range("A47").Value = 46070
$ws.Range("D47").Value = 159.94
$ws.Range("E47").Value = 151.08000000000001
$ws.Range("F47").Value = 161.08000000000001

$ws.Range("A48").Value = 46070
$ws.Range("D48").Value = 159.63
$ws.Range("E48").Value = 151.06
$ws.Range("F48").Value = 161.06

$ws.Range("A49").Value = 46067
$ws.Range("D49").Value = 160.02000000000001
$ws.Range("E49").Value = 151.36000000000001
$ws.Range("F49").Value = 161.36000000000001

$ws.Range("A50").Value = 46067
$ws.Range("D50").Value = 159.69999999999999
$ws.Range("E50").Value = 151.33000000000001
$ws.Range("F50").Value = 161.33000000000001

$ws.Range("A54").Value = 46070
$ws.Range("D54").Value = 173.59
$ws.Range("E54").Value = 164.53
$ws.Range("F54").Value = 174.53

$ws.Range("A55").Value = 46070
$ws.Range("D55").Value = 162.96
$ws.Range("E55").Value = 163.22
$ws.Range("F55").Value = 173.22

$ws.Range("A56").Value = 46070
$ws.Range("D56").Value = 162.72999999999999

$ws.Range("A57").Value = 46070
$ws.Range("D57").Value = 163.59
$ws.Range("E57").Value = 157.63999999999999

$ws.Range("A58").Value = 46070
$ws.Range("D58").Value = 159.36000000000001
$ws.Range("E58").Value = 153.54
$ws.Range("F58").Value = 163.54

$ws.Range("A59").Value = 46070
$ws.Range("D59").Value = 166.5
$ws.Range("E59").Value = 163.05000000000001

$ws.Range("A60").Value = 46067
$ws.Range("D60").Value = 174.45
$ws.Range("E60").Value = 165.33
$ws.Range("F60").Value = 175.33

$ws.Range("A61").Value = 46067
$ws.Range("D61").Value = 163.82
$ws.Range("E61").Value = 163.75
$ws.Range("F61").Value = 173.75

$ws.Range("A62").Value = 46067
$ws.Range("D62").Value = 163.58000000000001

$ws.Range("A63").Value = 46067
$ws.Range("D63").Value = 164.4
$ws.Range("E63").Value = 158.16999999999999

$ws.Range("A64").Value = 46067
$ws.Range("D64").Value = 160.16999999999999
$ws.Range("E64").Value = 154.07
$ws.Range("F64").Value = 164.07

$ws.Range("A65").Value = 46067
$ws.Range("D65").Value = 167.3
$ws.Range("E65").Value = 163.80000000000001
